# Atualização na modelagem do BD e Backlog
#
# The "Missão/ Valores/ Visão" backlog item (row 13) is folded into the
# "Pagina quem somos nós" item right above it, so the standalone row is
# removed and the remaining Website items are renumbered / reclassified.
# A couple of other classification (RF/RNF) fixes ride along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Merge "Missão/ Valores/ Visão" into "Pagina quem somos nós" and drop
#    the now-redundant row (row 13) - everything below shifts up by one.
$ws.Range("G12").Value = "Pagina quem somos nós(Missão/ Valores/ Visão)"
$ws.Rows("13:13").Delete()

# 2) Website section is now entirely "RF" (Requisito Funcional).
$ws.Range("H9").Value = "RF"
$ws.Range("H11").Value = "RF"
$ws.Range("H12").Value = "RF"

# 3) Renumber the two Website rows that shifted up past the removed item.
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 8

# 4) Classification fixes elsewhere in the backlog.
$ws.Range("H20").Value = "RNF"
$ws.Range("H24").Value = "RF"

# 5) Update the view so the selection matches where the user left off.
$ws.Range("C9").Select()
